$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1049893408696204
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0.1618706169496776
$ws.Range("C5").Value = 0.00000000000000003582521641092274
$ws.Range("C6").Value = 0.7325499012809557
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0.0004661450251538002
$ws.Range("C10").Value = 0.00000000000000004155213597677616
$ws.Range("C11").Value = 0.0001239958745925973
$ws.Range("C12").Value = 0.000000000000000002788770720769446
$ws.Range("C13").Value = 0.000000000000000001958378521762924
